$d = $word.ActiveDocument

# Replace the hashcode in the EObjectImpl proxy reference.
$d.Content.Find.Execute("1e5bc4e7", $true, $false, $false, $false, $false, $true, 1, $false, "3c72031c", 2)

# Update stack trace line numbers (M2DocEvaluator.java).
$d.Content.Find.Execute("M2DocEvaluator.java:540)", $true, $false, $false, $false, $false, $true, 1, $false, "M2DocEvaluator.java:543)", 2)
$d.Content.Find.Execute("M2DocEvaluator.java:1038)", $true, $false, $false, $false, $false, $true, 1, $false, "M2DocEvaluator.java:1084)", 2)
$d.Content.Find.Execute("M2DocEvaluator.java:1254)", $true, $false, $false, $false, $false, $true, 1, $false, "M2DocEvaluator.java:1300)", 2)
$d.Content.Find.Execute("M2DocEvaluator.java:275)", $true, $false, $false, $false, $false, $true, 1, $false, "M2DocEvaluator.java:278)", 2)
$d.Content.Find.Execute("M2DocEvaluator.java:264)", $true, $false, $false, $false, $false, $true, 1, $false, "M2DocEvaluator.java:267)", 2)

# Update stack trace line numbers (M2DocUtils.java).
$d.Content.Find.Execute("M2DocUtils.java:712)", $true, $false, $false, $false, $false, $true, 1, $false, "M2DocUtils.java:694)", 2)

# Update stack trace line numbers (AbstractTemplatesTestSuite.java).
$d.Content.Find.Execute("AbstractTemplatesTestSuite.java:459)", $true, $false, $false, $false, $false, $true, 1, $false, "AbstractTemplatesTestSuite.java:476)", 2)
$d.Content.Find.Execute("AbstractTemplatesTestSuite.java:369)", $true, $false, $false, $false, $false, $true, 1, $false, "AbstractTemplatesTestSuite.java:385)", 2)

# Update the generated method accessor number.
$d.Content.Find.Execute("GeneratedMethodAccessor107", $true, $false, $false, $false, $false, $true, 1, $false, "GeneratedMethodAccessor111", 2)
